$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33,8).Value = 871.4091
$ws.Cells.Item(33,9).Value = 692.8333
$ws.Cells.Item(33,10).Value = 1675
$ws.Cells.Item(33,11).Value = 692.8333
$ws.Cells.Item(33,12).Value = 1675
$ws.Cells.Item(33,13).Value = -463.8333
$ws.Cells.Item(33,14).Value = -2133

$ws.Cells.Item(49,8).Value = 1034
$ws.Cells.Item(49,9).Value = 1658.5
$ws.Cells.Item(49,10).Value = 409.5
$ws.Cells.Item(49,11).Value = 4975.5
$ws.Cells.Item(49,12).Value = 1228.5
$ws.Cells.Item(49,13).Value = -4839.5
$ws.Cells.Item(49,14).Value = -1500.5

$ws.Cells.Item(125,8).Value = 71428840
$ws.Cells.Item(125,9).Value = 316.4
$ws.Cells.Item(125,10).Value = 250000140
$ws.Cells.Item(125,11).Value = 2847.6
$ws.Cells.Item(125,12).Value = 2250001260
$ws.Cells.Item(125,13).Value = -387.5999999999999
$ws.Cells.Item(125,14).Value = -2250006180

$ws.Cells.Item(132,8).Value = 1977.9318
$ws.Cells.Item(132,9).Value = 1135.4
$ws.Cells.Item(132,10).Value = 5254.4443
$ws.Cells.Item(132,11).Value = 3406.2
$ws.Cells.Item(132,12).Value = 15763.3329
$ws.Cells.Item(132,13).Value = -876.2000000000003
$ws.Cells.Item(132,14).Value = -20823.3329

$ws.Cells.Item(138,8).Value = 1547.4572
$ws.Cells.Item(138,9).Value = 851.95746
$ws.Cells.Item(138,10).Value = 2968.6956
$ws.Cells.Item(138,11).Value = 2555.87238
$ws.Cells.Item(138,12).Value = 8906.086800000001
$ws.Cells.Item(138,13).Value = 2584.12762
$ws.Cells.Item(138,14).Value = -19186.0868

$ws.Cells.Item(141,8).Value = 706.3488
$ws.Cells.Item(141,9).Value = 601.619
$ws.Cells.Item(141,10).Value = 5105
$ws.Cells.Item(141,11).Value = 1804.857
$ws.Cells.Item(141,12).Value = 15315
$ws.Cells.Item(141,13).Value = 3375.143
$ws.Cells.Item(141,14).Value = -25675

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74,8).Value = 205450.19
$ws.Cells.Item(74,9).Value = 323605.25
$ws.Cells.Item(74,10).Value = 46197.695
$ws.Cells.Item(74,11).Value = 323605.25
$ws.Cells.Item(74,12).Value = 46197.695
$ws.Cells.Item(74,13).Value = -322731.25
$ws.Cells.Item(74,14).Value = -47945.695

$ws.Cells.Item(77,8).Value = 205450.19
$ws.Cells.Item(77,9).Value = 323605.25
$ws.Cells.Item(77,10).Value = 46197.695
$ws.Cells.Item(77,11).Value = 1618026.25
$ws.Cells.Item(77,12).Value = 230988.475
$ws.Cells.Item(77,13).Value = -1613658.25
$ws.Cells.Item(77,14).Value = -239724.475

$ws.Cells.Item(111,8).Value = 0
$ws.Cells.Item(111,9).Value = 0
$ws.Cells.Item(111,10).Value = 0
$ws.Cells.Item(111,11).Value = 0
$ws.Cells.Item(111,12).Value = 0
$ws.Cells.Item(111,14).ClearContents()

$ws.Cells.Item(122,8).Value = 4339.04
$ws.Cells.Item(122,9).Value = 4436.778
$ws.Cells.Item(122,10).Value = 4087.7144
$ws.Cells.Item(122,11).Value = 13310.334
$ws.Cells.Item(122,12).Value = 12263.1432
$ws.Cells.Item(122,13).Value = -10860.334
$ws.Cells.Item(122,14).Value = -17163.1432

$ws.Cells.Item(133,8).Value = 30695.25
$ws.Cells.Item(133,9).Value = 0
$ws.Cells.Item(133,10).Value = 30695.25
$ws.Cells.Item(133,11).Value = 0
$ws.Cells.Item(133,12).Value = 30695.25
$ws.Cells.Item(133,14).Value = -35755.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99,8).Value = 9239040
$ws.Cells.Item(99,9).Value = 3215093.8
$ws.Cells.Item(99,10).Value = 33334826
$ws.Cells.Item(99,11).Value = 3215093.8
$ws.Cells.Item(99,12).Value = 33334826
$ws.Cells.Item(99,13).Value = -3213595.8
$ws.Cells.Item(99,14).Value = -33337822

$ws.Cells.Item(105,8).Value = 1738.55
$ws.Cells.Item(105,9).Value = 1600.7142
$ws.Cells.Item(105,10).Value = 2060.1667
$ws.Cells.Item(105,11).Value = 1600.7142
$ws.Cells.Item(105,12).Value = 2060.1667
$ws.Cells.Item(105,13).Value = 146.2858000000001
$ws.Cells.Item(105,14).Value = -5554.1667

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(17,8).Value = 4752
$ws.Cells.Item(17,9).Value = 10008
$ws.Cells.Item(17,10).Value = 3000
$ws.Cells.Item(17,11).Value = 10008
$ws.Cells.Item(17,12).Value = 3000
$ws.Cells.Item(17,13).Value = -9834
$ws.Cells.Item(17,14).Value = -3348

$ws.Cells.Item(31,8).Value = 3255.818
$ws.Cells.Item(31,9).Value = 1891.6552
$ws.Cells.Item(31,10).Value = 5893.2
$ws.Cells.Item(31,11).Value = 1891.6552
$ws.Cells.Item(31,12).Value = 5893.2
$ws.Cells.Item(31,13).Value = -1596.6552
$ws.Cells.Item(31,14).Value = -6483.2

$ws.Cells.Item(34,8).Value = 3255.818
$ws.Cells.Item(34,9).Value = 1891.6552
$ws.Cells.Item(34,10).Value = 5893.2
$ws.Cells.Item(34,11).Value = 1891.6552
$ws.Cells.Item(34,12).Value = 5893.2
$ws.Cells.Item(34,13).Value = -1689.6552
$ws.Cells.Item(34,14).Value = -6297.2

$ws.Cells.Item(58,8).Value = 2503.7285
$ws.Cells.Item(58,9).Value = 2731.0625
$ws.Cells.Item(58,10).Value = 2007.7273
$ws.Cells.Item(58,11).Value = 2731.0625
$ws.Cells.Item(58,12).Value = 2007.7273
$ws.Cells.Item(58,13).Value = -2528.0625
$ws.Cells.Item(58,14).Value = -2413.7273

$ws.Cells.Item(64,8).Value = 0
$ws.Cells.Item(64,9).Value = 0
$ws.Cells.Item(64,10).Value = 0
$ws.Cells.Item(64,11).Value = 0
$ws.Cells.Item(64,12).Value = 0
$ws.Cells.Item(64,14).ClearContents()

$ws.Cells.Item(67,8).Value = 0
$ws.Cells.Item(67,9).Value = 0
$ws.Cells.Item(67,10).Value = 0
$ws.Cells.Item(67,11).Value = 0
$ws.Cells.Item(67,12).Value = 0
$ws.Cells.Item(67,14).ClearContents()

$ws.Cells.Item(99,8).Value = 64493.062
$ws.Cells.Item(99,9).Value = 92789.17999999999
$ws.Cells.Item(99,10).Value = 2241.6
$ws.Cells.Item(99,11).Value = 92789.17999999999
$ws.Cells.Item(99,12).Value = 2241.6
$ws.Cells.Item(99,13).Value = -91291.17999999999
$ws.Cells.Item(99,14).Value = -5237.6

$ws.Cells.Item(126,8).Value = 64493.062
$ws.Cells.Item(126,9).Value = 92789.17999999999
$ws.Cells.Item(126,10).Value = 2241.6
$ws.Cells.Item(126,11).Value = 278367.54
$ws.Cells.Item(126,12).Value = 6724.799999999999
$ws.Cells.Item(126,13).Value = -275897.54
$ws.Cells.Item(126,14).Value = -11664.8

$ws.Cells.Item(134,8).Value = 1443.8235
$ws.Cells.Item(134,9).Value = 987.3019
$ws.Cells.Item(134,10).Value = 2199.9375
$ws.Cells.Item(134,11).Value = 2961.9057
$ws.Cells.Item(134,12).Value = 6599.8125
$ws.Cells.Item(134,13).Value = -426.9057000000003
$ws.Cells.Item(134,14).Value = -11669.8125

$ws.Cells.Item(136,8).Value = 2503.7285
$ws.Cells.Item(136,9).Value = 2731.0625
$ws.Cells.Item(136,10).Value = 2007.7273
$ws.Cells.Item(136,11).Value = 8193.1875
$ws.Cells.Item(136,12).Value = 6023.1819
$ws.Cells.Item(136,13).Value = -5643.1875
$ws.Cells.Item(136,14).Value = -11123.1819

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5,8).Value = 558.8570999999999
$ws.Cells.Item(5,9).Value = 493.92
$ws.Cells.Item(5,10).Value = 1100
$ws.Cells.Item(5,11).Value = 1481.76
$ws.Cells.Item(5,12).Value = 3300
$ws.Cells.Item(5,13).Value = -1369.76
$ws.Cells.Item(5,14).Value = -3524

$ws.Cells.Item(74,8).Value = 2064.4443
$ws.Cells.Item(74,9).Value = 925
$ws.Cells.Item(74,10).Value = 2976
$ws.Cells.Item(74,11).Value = 2775
$ws.Cells.Item(74,12).Value = 8928
$ws.Cells.Item(74,13).Value = -1714
$ws.Cells.Item(74,14).Value = -11050

$ws.Cells.Item(77,8).Value = 2064.4443
$ws.Cells.Item(77,9).Value = 925
$ws.Cells.Item(77,10).Value = 2976
$ws.Cells.Item(77,11).Value = 8325
$ws.Cells.Item(77,12).Value = 26784
$ws.Cells.Item(77,13).Value = -3021
$ws.Cells.Item(77,14).Value = -37392

$ws.Cells.Item(122,8).Value = 990.381
$ws.Cells.Item(122,9).Value = 562.25
$ws.Cells.Item(122,10).Value = 1253.8462
$ws.Cells.Item(122,11).Value = 5060.25
$ws.Cells.Item(122,12).Value = 11284.6158
$ws.Cells.Item(122,13).Value = -2610.25
$ws.Cells.Item(122,14).Value = -16184.6158

$ws.Cells.Item(124,8).Value = 2316
$ws.Cells.Item(124,9).Value = 280
$ws.Cells.Item(124,10).Value = 2825
$ws.Cells.Item(124,11).Value = 840
$ws.Cells.Item(124,12).Value = 8475
$ws.Cells.Item(124,13).Value = 4070
$ws.Cells.Item(124,14).Value = -18295

$ws.Cells.Item(125,8).Value = 3336.6667
$ws.Cells.Item(125,9).Value = 530
$ws.Cells.Item(125,10).Value = 3898
$ws.Cells.Item(125,11).Value = 1590
$ws.Cells.Item(125,12).Value = 11694
$ws.Cells.Item(125,13).Value = 3330
$ws.Cells.Item(125,14).Value = -21534

$ws.Cells.Item(135,8).Value = 558.8570999999999
$ws.Cells.Item(135,9).Value = 493.92
$ws.Cells.Item(135,10).Value = 1100
$ws.Cells.Item(135,11).Value = 4445.28
$ws.Cells.Item(135,12).Value = 9900
$ws.Cells.Item(135,13).Value = -1910.28
$ws.Cells.Item(135,14).Value = -14970

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102,8).Value = 4290.9375
$ws.Cells.Item(102,9).Value = 1939.6428
$ws.Cells.Item(102,10).Value = 20750
$ws.Cells.Item(102,11).Value = 1939.6428
$ws.Cells.Item(102,12).Value = 20750
$ws.Cells.Item(102,13).Value = -317.6428000000001
$ws.Cells.Item(102,14).Value = -23994

$ws.Cells.Item(141,8).Value = 30588.6
$ws.Cells.Item(141,9).Value = 14555
$ws.Cells.Item(141,10).Value = 34597
$ws.Cells.Item(141,11).Value = 14555
$ws.Cells.Item(141,12).Value = 34597
$ws.Cells.Item(141,13).Value = -9375
$ws.Cells.Item(141,14).Value = -44957

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7,8).Value = 2155.4856
$ws.Cells.Item(7,9).Value = 1872
$ws.Cells.Item(7,10).Value = 2580.7144
$ws.Cells.Item(7,11).Value = 1872
$ws.Cells.Item(7,12).Value = 2580.7144
$ws.Cells.Item(7,13).Value = -1760
$ws.Cells.Item(7,14).Value = -2804.7144

$ws.Cells.Item(82,8).Value = 2900
$ws.Cells.Item(82,9).Value = 3137.5
$ws.Cells.Item(82,10).Value = 2710
$ws.Cells.Item(82,11).Value = 3137.5
$ws.Cells.Item(82,12).Value = 2710
$ws.Cells.Item(82,13).Value = -2776.5
$ws.Cells.Item(82,14).Value = -3432

$ws.Cells.Item(85,8).Value = 2900
$ws.Cells.Item(85,9).Value = 3137.5
$ws.Cells.Item(85,10).Value = 2710
$ws.Cells.Item(85,11).Value = 3137.5
$ws.Cells.Item(85,12).Value = 2710
$ws.Cells.Item(85,13).Value = -1889.5
$ws.Cells.Item(85,14).Value = -5206

$ws.Cells.Item(126,8).Value = 2155.4856
$ws.Cells.Item(126,9).Value = 1872
$ws.Cells.Item(126,10).Value = 2580.7144
$ws.Cells.Item(126,11).Value = 5616
$ws.Cells.Item(126,12).Value = 7742.1432
$ws.Cells.Item(126,13).Value = -3146
$ws.Cells.Item(126,14).Value = -12682.1432

$ws.Cells.Item(133,8).Value = 18172.445
$ws.Cells.Item(133,9).Value = 0
$ws.Cells.Item(133,10).Value = 18172.445
$ws.Cells.Item(133,11).Value = 0
$ws.Cells.Item(133,12).Value = 18172.445
$ws.Cells.Item(133,14).Value = -23232.445

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(80,8).Value = 0
$ws.Cells.Item(80,9).Value = 0
$ws.Cells.Item(80,10).Value = 0
$ws.Cells.Item(80,11).Value = 0
$ws.Cells.Item(80,12).Value = 0
$ws.Cells.Item(80,14).ClearContents()

$ws.Cells.Item(83,8).Value = 0
$ws.Cells.Item(83,9).Value = 0
$ws.Cells.Item(83,10).Value = 0
$ws.Cells.Item(83,11).Value = 0
$ws.Cells.Item(83,12).Value = 0
$ws.Cells.Item(83,14).ClearContents()

$ws.Cells.Item(136,8).Value = 11757419
$ws.Cells.Item(136,9).Value = 17562694
$ws.Cells.Item(136,10).Value = 347048.66
$ws.Cells.Item(136,11).Value = 52688082
$ws.Cells.Item(136,12).Value = 1041145.98
$ws.Cells.Item(136,13).Value = -52685532
$ws.Cells.Item(136,14).Value = -1046245.98
